# Apply license-upgrade edits to the item_license workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "shop_item" to "item_license"
$ws.Name = "item_license"

# Update the "value" column (D) for rows 3-10 to reflect new license prices.
$ws.Range("D3").Value = 2000
$ws.Range("D4").Value = 2000
$ws.Range("D5").Value = 3000
$ws.Range("D6").Value = 3000
$ws.Range("D7").Value = 3000
$ws.Range("D8").Value = 4000
$ws.Range("D9").Value = 4000
$ws.Range("D10").Value = 4000

# Move the active selection to D10, matching the saved view state.
$ws.Range("D10").Select()

# Update the window size/position to match the saved workbook view.
$win = $excel.Windows.Item(1)
$win.Left = 19200
$win.Top = 0
$win.Width = 19200
$win.Height = 10845
